$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 1 de Mayo de 2020 a las 11:22"

# Row 16 - Belgica
$ws.Range("B16").Value = 49032
$ws.Range("C16").Value = 513
$ws.Range("D16").Value = 11892
$ws.Range("E16").Value = 29437
$ws.Range("F16").Value = 740
$ws.Range("G16").Value = 109
$ws.Range("H16").Value = 7703

# Row 40 - Indonesia
$ws.Range("B40").Value = 10551
$ws.Range("C40").Value = 433
$ws.Range("D40").Value = 1591
$ws.Range("E40").Value = 8160
$ws.Range("G40").Value = 8
$ws.Range("H40").Value = 800

# Row 51 - Malasia
$ws.Range("B51").Value = 6071
$ws.Range("C51").Value = 69
$ws.Range("D51").Value = 4210
$ws.Range("E51").Value = 1758
$ws.Range("F51").Value = 37
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 103

# Row 59 - Moldavia
$ws.Range("D59").Value = 1272
$ws.Range("E59").Value = 2506
$ws.Range("G59").Value = 3
$ws.Range("H59").Value = 119

# Row 66 - Oman
$ws.Range("B66").Value = 2447
$ws.Range("C66").Value = 99
$ws.Range("E66").Value = 1941

# Row 103 - Sri Lanka
$ws.Range("B103").Value = 668
$ws.Range("C103").Value = 5
$ws.Range("E103").Value = 504

# Row 142 - Etiopia
$ws.Range("B142").Value = 133
$ws.Range("C142").Value = 2
$ws.Range("D142").Value = 66
$ws.Range("E142").Value = 64

$wb.Save()
